$d = $word.ActiveDocument

# 1. Title: "Home Assignment 1" -> "HOME ASSIGNMENT 1" and font size 26 -> 32 (both pPr/rPr and run rPr)
$para = $d.Paragraphs(1)
$para.Range.Font.Size = 16
$para.Range.ParagraphFormat.Style.Font.Size = 16

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Home Assignment 1", $true, $false, $false, $false, $false, `
              $true, 1, $false, "HOME ASSIGNMENT 1", 2)

# 2. "Project Name" -> "PROJECT NAME"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("Project Name", $true, $false, $false, $false, $false, `
               $true, 1, $false, "PROJECT NAME", 2)

# 3. Style change: add semiHidden to DefaultParagraphFont style
$style = $d.Styles("Default Paragraph Font")
$style.SemiHidden = $true

$d.Save()
